$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 254.77777
$ws.Range("I2").Value = 98.71429000000001
$ws.Range("K2").Value = 98.71429000000001
$ws.Range("M2").Value = 14.28570999999999
$ws.Range("H8").Value = 372
$ws.Range("I8").Value = 182.5
$ws.Range("J8").Value = 751
$ws.Range("K8").Value = 547.5
$ws.Range("L8").Value = 2253
$ws.Range("M8").Value = -408.5
$ws.Range("N8").Value = -2531
$ws.Range("H12").Value = 1083.4445
$ws.Range("I12").Value = 686.1429000000001
$ws.Range("K12").Value = 686.1429000000001
$ws.Range("M12").Value = -516.1429000000001
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 9900
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H18").Value = 325
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H19").Value = 1144.8667
$ws.Range("I19").Value = 1326.6
$ws.Range("J19").Value = 1054
$ws.Range("K19").Value = 1326.6
$ws.Range("L19").Value = 1054
$ws.Range("M19").Value = -1151.6
$ws.Range("N19").Value = -1404
$ws.Range("H26").Value = 59999
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H33").Value = 469.53333
$ws.Range("J33").Value = 1490.6666
$ws.Range("L33").Value = 1490.6666
$ws.Range("N33").Value = -1948.6666
$ws.Range("H34").Value = 5393.9
$ws.Range("I34").Value = 5393.9
$ws.Range("K34").Value = 5393.9
$ws.Range("M34").Value = -5190.9
$ws.Range("H36").Value = 5393.9
$ws.Range("I36").Value = 5393.9
$ws.Range("K36").Value = 5393.9
$ws.Range("M36").Value = -4678.9
$ws.Range("H40").Value = 3869.5217
$ws.Range("I40").Value = 2999.6667
$ws.Range("K40").Value = 2999.6667
$ws.Range("M40").Value = -2824.6667
$ws.Range("H41").Value = 5575.4287
$ws.Range("I41").Value = 5071.6665
$ws.Range("J41").Value = 8598
$ws.Range("K41").Value = 5071.6665
$ws.Range("L41").Value = 8598
$ws.Range("M41").Value = -4631.6665
$ws.Range("N41").Value = -9478
$ws.Range("H42").Value = 116.57143
$ws.Range("J42").Value = 86
$ws.Range("L42").Value = 258
$ws.Range("N42").Value = -718
$ws.Range("H43").Value = 1390.8182
$ws.Range("I43").Value = 1409.9
$ws.Range("K43").Value = 1409.9
$ws.Range("M43").Value = -1340.9
$ws.Range("H44").Value = 40012.5
$ws.Range("J44").Value = 40012.5
$ws.Range("L44").Value = 40012.5
$ws.Range("N44").Value = -40936.5
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H51").Value = 3397.25
$ws.Range("I51").Value = 3714.842
$ws.Range("J51").Value = 2190.4
$ws.Range("K51").Value = 3714.842
$ws.Range("L51").Value = 2190.4
$ws.Range("M51").Value = -3230.842
$ws.Range("N51").Value = -3158.4
$ws.Range("H53").Value = 1010.65
$ws.Range("I53").Value = 976.26666
$ws.Range("K53").Value = 976.26666
$ws.Range("M53").Value = -339.26666
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H55").Value = 362
$ws.Range("I55").Value = 359.8889
$ws.Range("J55").Value = 368.33334
$ws.Range("K55").Value = 359.8889
$ws.Range("L55").Value = 368.33334
$ws.Range("M55").Value = -145.8889
$ws.Range("N55").Value = -796.33334
$ws.Range("H58").Value = 1305.6
$ws.Range("I58").Value = 842.6667
$ws.Range("K58").Value = 2528.0001
$ws.Range("M58").Value = -2378.0001
$ws.Range("H61").Value = 142.85715
$ws.Range("I61").Value = 142.85715
$ws.Range("K61").Value = 428.57145
$ws.Range("M61").Value = -256.57145
$ws.Range("H100").Value = 2748.1052
$ws.Range("I100").Value = 1421.4
$ws.Range("K100").Value = 1421.4
$ws.Range("M100").Value = -880.4000000000001
$ws.Range("H106").Value = 4168796.5
$ws.Range("I106").Value = 4168796.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4168796.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4168165.5
$ws.Range("N106").ClearContents()
$ws.Range("H121").Value = 2964.7778
$ws.Range("J121").Value = 2964.7778
$ws.Range("L121").Value = 8894.3334
$ws.Range("N121").Value = -12388.3334
$ws.Range("H126").Value = 150990
$ws.Range("J126").Value = 150990
$ws.Range("L126").Value = 150990
$ws.Range("N126").Value = -160870
$ws.Range("H128").Value = 84981.25
$ws.Range("J128").Value = 84981.25
$ws.Range("L128").Value = 84981.25
$ws.Range("N128").Value = -94941.25
$ws.Range("H130").Value = 116998
$ws.Range("J130").Value = 116998
$ws.Range("L130").Value = 116998
$ws.Range("N130").Value = -127038
$ws.Range("H137").Value = 9264667
$ws.Range("J137").Value = 9264667
$ws.Range("L137").Value = 27794001
$ws.Range("N137").Value = -27799101
$ws.Range("H138").Value = 4486.4634
$ws.Range("J138").Value = 4745.6763
$ws.Range("L138").Value = 14237.0289
$ws.Range("N138").Value = -24517.0289

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3962.2424
$ws.Range("I122").Value = 3382.9048
$ws.Range("K122").Value = 10148.7144
$ws.Range("M122").Value = -7698.714399999999
$ws.Range("H132").Value = 3780.238
$ws.Range("I132").Value = 3780.238
$ws.Range("K132").Value = 11340.714
$ws.Range("M132").Value = -8810.714

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 812.95
$ws.Range("I94").Value = 961.1667
$ws.Range("K94").Value = 961.1667
$ws.Range("M94").Value = -510.1667
$ws.Range("H117").Value = 144000
$ws.Range("J117").Value = 144000
$ws.Range("L117").Value = 144000
$ws.Range("N117").Value = -153178
$ws.Range("H141").Value = 126480
$ws.Range("J141").Value = 126480
$ws.Range("L141").Value = 126480
$ws.Range("N141").Value = -136840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 3010
$ws.Range("I14").Value = 3010
$ws.Range("K14").Value = 3010
$ws.Range("M14").Value = -2840
$ws.Range("H16").Value = 2090.8
$ws.Range("I16").Value = 2043.5
$ws.Range("K16").Value = 2043.5
$ws.Range("M16").Value = -1756.5
$ws.Range("H21").Value = 4998.3335
$ws.Range("I21").Value = 4998.3335
$ws.Range("K21").Value = 4998.3335
$ws.Range("M21").Value = -4763.3335
$ws.Range("H31").Value = 7507.7144
$ws.Range("I31").Value = 3896.125
$ws.Range("K31").Value = 3896.125
$ws.Range("M31").Value = -3601.125
$ws.Range("H34").Value = 7507.7144
$ws.Range("I34").Value = 3896.125
$ws.Range("K34").Value = 3896.125
$ws.Range("M34").Value = -3694.125
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 42499.25
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4264
$ws.Range("H61").Value = 42499.25
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4652
$ws.Range("H113").Value = 2090.8
$ws.Range("I113").Value = 2043.5
$ws.Range("K113").Value = 2043.5
$ws.Range("M113").Value = 126.5
$ws.Range("H114").Value = 79978.5
$ws.Range("J114").Value = 79978.5
$ws.Range("L114").Value = 79978.5
$ws.Range("N114").Value = -88656.5
$ws.Range("H133").Value = 29887.5
$ws.Range("J133").Value = 29888
$ws.Range("L133").Value = 29888
$ws.Range("N133").Value = -34948
$ws.Range("H134").Value = 3077.125
$ws.Range("I134").Value = 3077.125
$ws.Range("K134").Value = 9231.375
$ws.Range("M134").Value = -6696.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7578.4
$ws.Range("I3").Value = 7578.4
$ws.Range("K3").Value = 22735.2
$ws.Range("M3").Value = -22623.2
$ws.Range("H35").Value = 2251.2
$ws.Range("I35").Value = 2251.2
$ws.Range("K35").Value = 6753.599999999999
$ws.Range("M35").Value = -6465.599999999999
$ws.Range("H36").Value = 1850
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 700
$ws.Range("K36").Value = 9000
$ws.Range("L36").Value = 2100
$ws.Range("M36").Value = -8831
$ws.Range("N36").Value = -2438
$ws.Range("H101").Value = 1000
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -7868
$ws.Range("H131").Value = 1669.2413
$ws.Range("J131").Value = 1819.3182
$ws.Range("L131").Value = 5457.9546
$ws.Range("N131").Value = -15537.9546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 321.85
$ws.Range("I2").Value = 251.13333
$ws.Range("K2").Value = 251.13333
$ws.Range("M2").Value = -138.13333
$ws.Range("H97").Value = 661.46155
$ws.Range("I97").Value = 457.2
$ws.Range("K97").Value = 457.2
$ws.Range("M97").Value = 38.80000000000001
$ws.Range("H107").Value = 999.6667
$ws.Range("I107").Value = 999.6667
$ws.Range("K107").Value = 999.6667
$ws.Range("M107").Value = 920.3333
$ws.Range("H122").Value = 3064.5833
$ws.Range("I122").Value = 3078.5
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 9235.5
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -6785.5
$ws.Range("N122").Value = -13885
$ws.Range("H129").Value = 121960
$ws.Range("J129").Value = 121960
$ws.Range("L129").Value = 121960
$ws.Range("N129").Value = -131960
$ws.Range("H132").Value = 2731.9
$ws.Range("I132").Value = 2731.9
$ws.Range("K132").Value = 8195.700000000001
$ws.Range("M132").Value = -5665.700000000001
$ws.Range("H133").Value = 52500
$ws.Range("J133").Value = 52500
$ws.Range("L133").Value = 52500
$ws.Range("N133").Value = -62620
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 11502
$ws.Range("I3").Value = 8004
$ws.Range("K3").Value = 8004
$ws.Range("M3").Value = -7892
$ws.Range("H15").Value = 11502
$ws.Range("I15").Value = 8004
$ws.Range("K15").Value = 8004
$ws.Range("M15").Value = -7834
$ws.Range("H16").Value = 1139.3043
$ws.Range("I16").Value = 1190.619
$ws.Range("J16").Value = 600.5
$ws.Range("K16").Value = 1190.619
$ws.Range("L16").Value = 600.5
$ws.Range("M16").Value = -1020.619
$ws.Range("N16").Value = -940.5
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H24").Value = 14145.25
$ws.Range("I24").Value = 14145.25
$ws.Range("K24").Value = 14145.25
$ws.Range("M24").Value = -13802.25
$ws.Range("H25").Value = 12007
$ws.Range("I25").Value = 12007
$ws.Range("K25").Value = 12007
$ws.Range("M25").Value = -11777
$ws.Range("H46").Value = 3784
$ws.Range("I46").Value = 1849.5
$ws.Range("J46").Value = 3938.76
$ws.Range("K46").Value = 1849.5
$ws.Range("L46").Value = 3938.76
$ws.Range("M46").Value = -1661.5
$ws.Range("N46").Value = -4314.76
$ws.Range("H55").Value = 834.6579
$ws.Range("I55").Value = 753.6129
$ws.Range("J55").Value = 1193.5714
$ws.Range("K55").Value = 753.6129
$ws.Range("L55").Value = 1193.5714
$ws.Range("M55").Value = -580.6129
$ws.Range("N55").Value = -1539.5714
$ws.Range("H61").Value = 4664
$ws.Range("I61").Value = 1687.8889
$ws.Range("K61").Value = 1687.8889
$ws.Range("M61").Value = -1485.8889
$ws.Range("H68").Value = 2165.889
$ws.Range("J68").Value = 3003
$ws.Range("L68").Value = 3003
$ws.Range("N68").Value = -4501
$ws.Range("H71").Value = 2165.889
$ws.Range("J71").Value = 3003
$ws.Range("L71").Value = 15015
$ws.Range("N71").Value = -22503
$ws.Range("H113").Value = 4664
$ws.Range("I113").Value = 1687.8889
$ws.Range("K113").Value = 1687.8889
$ws.Range("M113").Value = 482.1111000000001
$ws.Range("H132").Value = 3299.6
$ws.Range("I132").Value = 3299.6
$ws.Range("K132").Value = 9898.799999999999
$ws.Range("M132").Value = -7368.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2509000
$ws.Range("I2").Value = 18000
$ws.Range("J2").Value = 5000000
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 5000000
$ws.Range("M2").Value = -17888
$ws.Range("N2").Value = -5000224
$ws.Range("H3").Value = 5501.5
$ws.Range("I3").Value = 5501.5
$ws.Range("K3").Value = 5501.5
$ws.Range("M3").Value = -5387.5
$ws.Range("H9").Value = 3000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 3000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3280
$ws.Range("H58").Value = 30085
$ws.Range("I58").Value = 30085
$ws.Range("K58").Value = 30085
$ws.Range("M58").Value = -29777
$ws.Range("H61").Value = 27421.2
$ws.Range("J61").Value = 39999.5
$ws.Range("L61").Value = 39999.5
$ws.Range("N61").Value = -40583.5
$ws.Range("H62").Value = 5087.5
$ws.Range("J62").Value = 6140
$ws.Range("L62").Value = 6140
$ws.Range("N62").Value = -7388
$ws.Range("H65").Value = 5087.5
$ws.Range("J65").Value = 6140
$ws.Range("L65").Value = 30700
$ws.Range("N65").Value = -36940
$ws.Range("H112").Value = 72453
$ws.Range("J112").Value = 72453
$ws.Range("L112").Value = 72453
$ws.Range("N112").Value = -75407
$ws.Range("H132").Value = 2435.037
$ws.Range("I132").Value = 2380.4783
$ws.Range("J132").Value = 2748.75
$ws.Range("K132").Value = 7141.4349
$ws.Range("L132").Value = 8246.25
$ws.Range("M132").Value = -4611.4349
$ws.Range("N132").Value = -13306.25
$ws.Range("H136").Value = 2725.7083
$ws.Range("I136").Value = 2072.3
$ws.Range("J136").Value = 5992.75
$ws.Range("K136").Value = 6216.900000000001
$ws.Range("L136").Value = 17978.25
$ws.Range("M136").Value = -3666.900000000001
$ws.Range("N136").Value = -23078.25
$ws.Range("H137").Value = 70715
$ws.Range("J137").Value = 70715
$ws.Range("L137").Value = 70715
$ws.Range("N137").Value = -80915
$ws.Range("H138").Value = 49429
$ws.Range("J138").Value = 49429
$ws.Range("N138").Value = -59709
